# Updating filtered feeds from workflow
# Appends one new row (58) to the "filtered feeds" sheet with a new link,
# its keywords and its title, mirroring the rows already present.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 58

$link = "https://www.genomeweb.com/companion-diagnostics/celebrating-10th-anniversary-its-pd-l1-cdx-agilent-sets-sights-next-gen"
$keywords = "CDx, companion diagnostic"
$title = "Celebrating the 10th Anniversary of its PD-L1 CDx, Agilent Sets Sights on Next-Gen Diagnostic Tools"

$ws.Cells.Item($newRow, 1).Value = $link
$ws.Cells.Item($newRow, 2).Value = $keywords
$ws.Cells.Item($newRow, 3).Value = $title

# Turn A<newRow> into a real hyperlink pointing at the link, same as every
# other row in column A.
$ws.Hyperlinks.Add($ws.Cells.Item($newRow, 1), $link) | Out-Null

# Match the link-cell style used by the rest of column A.
$ws.Cells.Item($newRow, 1).Style = $ws.Cells.Item($newRow - 1, 1).Style
